$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell text to reflect the reordered / renamed notes.
# A1 stays the same ("Finetune the third fold...")
$ws.Range("A2").Value = "Add TTA to predictions OK"
$ws.Range("A3").Value = "Add classifier of empty predictions"
$ws.Range("A4").Value = "Ensemble with resnet101 OK"
$ws.Range("A5").Value = "Remove few pixel images and retrain"

# Move the active selection from A6 to A3
$ws.Range("A3").Select()
